$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings that can look numeric (e.g. "248.64" or
# "1.001"); force text formatting on exactly the contiguous blocks of
# cells being rewritten so Excel does not silently convert the new values
# to numbers on assignment (each block set separately - a single
# multi-area Range(...) only honours NumberFormat on its first area).
$ws.Range("D2:D4").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8:D21").NumberFormat = "@"
$ws.Range("D23:D43").NumberFormat = "@"
$ws.Range("D45:D51").NumberFormat = "@"

$ws.Range("D2").Value = '30.395.07'
$ws.Range("E2").Value = '  +0.24%  '
$ws.Range("D3").Value = '1.938.15'
$ws.Range("E3").Value = '  +0.17%  '
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("E5").Value = '  +6.66%  '
$ws.Range("D6").Value = '248.64'
$ws.Range("E6").Value = '  -0.69%  '
$ws.Range("E7").Value = '  -0.06%  '
$ws.Range("D8").Value = '28.00'
$ws.Range("E8").Value = '  +2.27%  '
$ws.Range("D9").Value = '0.3206'
$ws.Range("E9").Value = '  -2.98%  '
$ws.Range("D10").Value = '0.07120'
$ws.Range("E10").Value = '  -1.67%  '
$ws.Range("D11").Value = '0.7866'
$ws.Range("E11").Value = '  -2.49%  '
$ws.Range("D12").Value = '0.08010'
$ws.Range("E12").Value = '  -0.80%  '
$ws.Range("D13").Value = '1.940.64'
$ws.Range("E13").Value = '  +0.36%  '
$ws.Range("D14").Value = '5.383'
$ws.Range("E14").Value = '  -1.76%  '
$ws.Range("D15").Value = '94.99'
$ws.Range("E15").Value = '  +0.71%  '
$ws.Range("D16").Value = '14.58'
$ws.Range("E16").Value = '  -3.46%  '
$ws.Range("D17").Value = '30.399.31'
$ws.Range("E17").Value = '  +0.27%  '
$ws.Range("D18").Value = '257.13'
$ws.Range("E18").Value = '  +2.68%  '
$ws.Range("D19").Value = '0.000008032'
$ws.Range("E19").Value = '  -1.84%  '
$ws.Range("D20").Value = '5.813'
$ws.Range("E20").Value = '  -0.28%  '
$ws.Range("D21").Value = '2.191.29'
$ws.Range("E21").Value = '  +0.29%  '
$ws.Range("E22").Value = '  -0.08%  '
$ws.Range("D23").Value = '1.001'
$ws.Range("E23").Value = '  +0.01%  '
$ws.Range("D24").Value = '6.789'
$ws.Range("E24").Value = '  -2.69%  '
$ws.Range("D25").Value = '9.612'
$ws.Range("E25").Value = '  -1.26%  '
$ws.Range("D26").Value = '164.88'
$ws.Range("E26").Value = '  +0.77%  '
$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").Value = '19.18'
$ws.Range("E27").Value = '  -0.39%  '
$ws.Range("B28").Value = 'Stellar'
$ws.Range("C28").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D28").Value = '0.1340'
$ws.Range("E28").Value = '  +1.67%  '
$ws.Range("D29").Value = '2.295'
$ws.Range("E29").Value = '  -3.18%  '
$ws.Range("D30").Value = '1.369'
$ws.Range("E30").Value = '  +1.47%  '
$ws.Range("D31").Value = '1.532'
$ws.Range("E31").Value = '  -2.13%  '
$ws.Range("D32").Value = '4.436'
$ws.Range("E32").Value = '  +0.72%  '
$ws.Range("D33").Value = '4.167'
$ws.Range("E33").Value = '  +0.46%  '
$ws.Range("D34").Value = '0.05205'
$ws.Range("E34").Value = '  +0.23%  '
$ws.Range("D35").Value = '1.282'
$ws.Range("E35").Value = '  +0.72%  '
$ws.Range("D36").Value = '0.7515'
$ws.Range("E36").Value = '  +1.11%  '
$ws.Range("D37").Value = '2.778'
$ws.Range("E37").Value = '  +1.13%  '
$ws.Range("D38").Value = '0.01973'
$ws.Range("E38").Value = '  +0.22%  '
$ws.Range("D39").Value = '2.806'
$ws.Range("E39").Value = '  -0.47%  '
$ws.Range("D40").Value = '78.12'
$ws.Range("E40").Value = '  -0.32%  '
$ws.Range("D41").Value = '6.464'
$ws.Range("E41").Value = '  +1.00%  '
$ws.Range("D42").Value = '0.4518'
$ws.Range("D43").Value = '1.985'
$ws.Range("E43").Value = '  -1.25%  '
$ws.Range("E44").Value = '  +0.08%  '
$ws.Range("D45").Value = '0.8383'
$ws.Range("E45").Value = '  -0.93%  '
$ws.Range("D46").Value = '101.72'
$ws.Range("E46").Value = '  +0.30%  '
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").Value = '9.823'
$ws.Range("E47").Value = '  +0.98%  '
$ws.Range("B48").Value = 'Aptos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D48").Value = '7.544'
$ws.Range("E48").Value = '  +1.51%  '
$ws.Range("D49").Value = '37.59'
$ws.Range("E49").Value = '  +2.83%  '
$ws.Range("D50").Value = '975.10'
$ws.Range("E50").Value = '  +10.71%  '
$ws.Range("B51").Value = 'NEARProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D51").Value = '1.505'
$ws.Range("E51").Value = '  +1.57%  '
